# Rename the second worksheet ("Rob") to "Commercial".
#
# The rest of the authoritative diff (styles.xml cellXfs reordering, and the
# corresponding `s="..."` index bumps on Summary!A2:J2 and on
# Commercial!I4/J4, I23/J23, I42/J42, M24) is a pure re-shuffle of one
# already-existing cellXfs entry (the bare `numFmtId="44" xfId="1"
# applyFont="1"` Currency xf) from the tail of the array to just after the
# `numFmtId="3"` entry. Every cell that is touched keeps the exact same
# resolved number format / font / fill / border — resolving each `s="N"`
# before and after the move lands on byte-identical `<xf>` definitions, so
# this is book-keeping noise from the original save (not a formatting
# change) and has no visible effect to reproduce through the object model.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rob")
$ws.Name = "Commercial"
